# Modify figure name for revised manuscript to eLife format.
#
# The manuscript figure built from this workbook was renamed for the
# eLife-formatted revision, so the three worksheets get more descriptive,
# reviewer-facing names. The per-lineage summary blocks at the bottom of
# each sheet also get trimmed down to just the columns actually used by
# the figure (A:E) - the leftover QUARTILE/AVERAGE/COUNTIF/ratio helper
# columns (F:K, F:I) are cleared out.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets -----------------------------------------------
$ws1 = $wb.Worksheets.Item("dividing_cell")
$ws1.Name = "growth-restored cell lineage"

$ws2 = $wb.Worksheets.Item("non-dividing_cell")
$ws2.Name = "growth-halted cell lineage"

$ws3 = $wb.Worksheets.Item("resistant_cell")
$ws3.Name = "non-deleted cell lineage"

# --- Clear the extra QUARTILE/AVERAGE/COUNTIF/ratio columns -----------
# Sheet 1 ("growth-restored cell lineage"): summary block rows 446:455
$ws1.Activate()
$ws1.Range("F446:K455").Select()
$ws1.Range("F446:K455").ClearContents()

# Sheet 2 ("growth-halted cell lineage"): summary block rows 745:754
$ws2.Activate()
$ws2.Range("F745:K754").Select()
$ws2.Range("F745:K754").ClearContents()

# Sheet 3 ("non-deleted cell lineage"): summary block rows 981:990
$ws3.Activate()
$ws3.Range("F981:I990").Select()
$ws3.Range("F981:I990").ClearContents()

# The saved workbook reopens on the third sheet.
$ws3.Activate()
